# forYaml.xlsx edit: populate the Sheet1 "YAML-as-cells" record for a
# simple WTG_onshore VariableRenewableOperator agent (used by the
# "simple future prices linear regression" commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: a lone numeric value in B1 (already bold/bordered via existing style)
$ws.Range("B1").Value = 0

# Row 2: name / WTG_onshore
$ws.Range("A2").Value = "name"
$ws.Range("B2").Value = "WTG_onshore"

# Row 3: parameters / {}
$ws.Range("A3").Value = "parameters"
$ws.Range("B3").Value = "{}"

# Row 4: type / VariableRenewableOperator
$ws.Range("A4").Value = "type"
$ws.Range("B4").Value = "VariableRenewableOperator"

# Row 5: InstalledPowerInMW / 300
$ws.Range("A5").Value = "InstalledPowerInMW"
$ws.Range("B5").Value = 300

# Row 6: technology / (blank value, cell stays present but empty)
$ws.Range("A6").Value = "technology"
$ws.Range("B6").Value = ""

# Reset selection back to A1 (originally parked at E9)
$ws.Range("A1").Select()
